$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("C1").Value = "archivo imagen"

# --- Data row (row 2) ---
$ws.Range("C2").Value = "imagen.jpg"
$ws.Range("E2").Value = "También llamada Gran Explosión, por su traducción en inglés."
$ws.Range("F2").Value = "densidad"
$ws.Range("H2").Value = "Stephen Hawking"
$zwsp = [char]0x200B
$ws.Range("I2").Value = "Stephen William Hawking (Oxford, 8 de enero de 1942-Cambridge, 14 de marzo de 2018)" + $zwsp + " fue un físico teórico, astrofísico, cosmólogo y divulgador científico británico."
$ws.Range("J2").Value = "singularidad"
$ws.Range("K2").Value = "Una singularidad gravitacional o espaciotemporal, de modo informal y desde un punto de vista físico, puede definirse como una zona del espacio-tiempo donde no se puede definir alguna magnitud física relacionada con los campos gravitatorios, tales como la curvatura, u otras."
$ws.Range("V2").Value = "Einstein"
$ws.Range("W2").Value = "Albert Einstein (Ulm, Imperio alemán, 14 de marzo de 1879-Princeton, Estados Unidos, 18 de abril de 1955) fue un físico alemán de origen judío, nacionalizado después suizo, austriaco y estadounidense. Se lo considera el científico más importante, conocido y popular del siglo XX.1"

# Wrap text for the newly long cells (matches style index used by existing wrapped cell G2)
$ws.Range("E2").WrapText = $true
$ws.Range("H2:W2").WrapText = $true

# Row height for row 2
$ws.Rows.Item(2).RowHeight = 105

# Column width: column W (23) becomes wider, split off from the former L:W group
$ws.Columns.Item(23).ColumnWidth = 41.083333333333336

# Update the active selection to G2
$ws.Range("G2").Select()
